$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "I was not particularly close..." -- em dash doubled
$ws.Range("B4").Value = "I was not particularly close to the Lord——we only met twice in total."

# Row 8: "After lunch..." -- Steward -> Butler
$ws.Range("B8").Value = "After lunch, I remained in my room, meditating, until Butler He came to inform me of the banquet’s time and location."

# Row 9: "At around..." -- time phrase changed; row shrinks from ht=51 to ht=34
$ws.Range("B9").Value = "At around 6.15 PM, I left my room just before the rain began and happened to run into Ming on the way."
$ws.Rows.Item(9).RowHeight = 34

# Row 12: wrap in green color tag + parentheses; row grows from ht=34 to ht=51
$ws.Range("B12").Value = " <color=#00CC00>(So, Kong was also someone who arrived early at the banquet and stayed the entire time.)</color>"
$ws.Rows.Item(12).RowHeight = 51

# Row 13: wrap in green color tag + parentheses
$ws.Range("B13").Value = " <color=#00CC00>(If that’s the case, it seems he didn’t have the opportunity to commit the crime.)</color>"

# Selection moves to B15
$ws.Range("B15").Select()
